$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - new track: "Utilities-Mine or System activation-13"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Utilities-Mine or System activation-13"
$ws.Range("C5").Value = 0.0020833333333333333
$ws.Range("C5").NumberFormat = "h:mm"
$ws.Range("D5").Value = "Future Weapons"
$ws.Range("E5").Value = "SoundMorph"
$ws.Range("E5").HorizontalAlignment = -4108
$ws.Range("F5").Value = "Start Hack, End Hack"
$ws.Range("G5").Value = "Sonniss.com - GDC 2015 - Game Audio Bundle"

# Row 6 - new track: "transition t04 soft 016"
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "transition t04 soft 016"
$ws.Range("C6").Value = 0.00069444444444444447
$ws.Range("C6").NumberFormat = "h:mm"
$ws.Range("D6").Value = "Transitions HD"
$ws.Range("E6").Value = "Alexander Kopeikin"
$ws.Range("F6").Value = "Spectral Jump Skill"
$ws.Range("G6").Value = "Sonniss.com - GDC 2015 - Game Audio Bundle"

# Update active selection to B8 (as reflected in the saved sheet view)
$ws.Range("B8").Select()
